$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.485.89'
$ws.Range('E2').Value = '  +3.13%  '
$ws.Range('D3').Value = '1.606.69'
$ws.Range('E3').Value = '  +2.85%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.75'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.12%  '
$ws.Range('E6').Value = '  +6.93%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '26.88'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +7.29%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '43.53'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.73%  '
$ws.Range('E10').Value = '  +2.78%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0602'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +2.80%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0911'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +1.66%  '
$ws.Range('D13').Value = '1.837.25'
$ws.Range('E13').Value = '  +2.84%  '
$ws.Range('D14').Value = '1.615.46'
$ws.Range('E14').Value = '  +3.34%  '
$ws.Range('D15').Value = '29.512.58'
$ws.Range('E15').Value = '  +3.20%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.535'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +4.12%  '
$ws.Range('E17').Value = '  +2.04%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '63.26'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +3.34%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '241.19'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +5.31%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.62'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +3.77%  '
$ws.Range('D21').Value = '0.0₃0691'
$ws.Range('E21').Value = '  +2.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.998'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.06%  '
$ws.Range('E23').Value = '  +2.61%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.20'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +2.36%  '
$ws.Range('E25').Value = '  +0.66%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '154.42'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +2.32%  '
$ws.Range('E27').Value = '  +5.12%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.26'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +3.39%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.38'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +2.49%  '
$ws.Range('E30').Value = '  -0.06%  '
$ws.Range('E31').Value = '  +2.52%  '
$ws.Range('E32').Value = '  +1.22%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.22'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +1.69%  '
$ws.Range('E34').Value = '  +4.24%  '
$ws.Range('D35').Value = '1.412.28'
$ws.Range('E35').Value = '  +1.86%  '
$ws.Range('E36').Value = '  +0.62%  '
$ws.Range('E37').Value = '  +3.60%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.84'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +5.51%  '
$ws.Range('E39').Value = '  +0.29%  '
$ws.Range('E40').Value = '  +2.72%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.538'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +3.63%  '
$ws.Range('E42').Value = '  +1.37%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0486'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +5.52%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.798'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +3.52%  '
$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.999'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.01%  '
$ws.Range('B46').Value = 'BitcoinSV'
$ws.Range('C46').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '52.70'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +22.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '65.62'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +2.62%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.28'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.00%  '
$ws.Range('D49').Value = '1.748.96'
$ws.Range('E49').Value = '  +3.02%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.861'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.00%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '86.68'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +1.93%  '
